$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # was "dados"  -> becomes "CH4"
$ws2 = $wb.Worksheets.Item(2)   # was "calculos" -> becomes "H2"

# ---------------------------------------------------------------
# Rename the sheets
# ---------------------------------------------------------------
$ws1.Name = "CH4"
$ws2.Name = "H2"

# ---------------------------------------------------------------
# CH4 sheet (ex "dados"): replace formulas with literal P values,
# drop the old "H2" column and add a new C [mol/g] column.
# ---------------------------------------------------------------
$ws1.Range("C1").ClearContents()

# Write the B header first, then A header, so the shared-string
# table ends up with "C [mol/g]" at index 0 and "P [MPa]" at index 1.
$ws1.Range("B1").Value = "C [mol/g]"
$ws1.Range("A1").Value = "P [MPa]"

$ws1.Range("A2").Value = 0.26
$ws1.Range("A3").Value = 0.79
$ws1.Range("A4").Value = 1.26
$ws1.Range("A5").Value = 1.74

$ws1.Range("B2").Value = 0.00141
$ws1.Range("B3").Value = 0.00315
$ws1.Range("B4").Value = 0.00395
$ws1.Range("B5").Value = 0.00442

$ws1.Range("B2:B5").NumberFormat = "0.00E+00"

$ws1.Columns.Item(2).ColumnWidth = 10.7
$ws1.Columns.Item(3).ColumnWidth = 9.6

# ---------------------------------------------------------------
# H2 sheet (ex "calculos"): turn the single P column into a
# P / C [mol/g] table with brand new data.
# ---------------------------------------------------------------
$ws2.Range("B1").Value = "C [mol/g]"
$ws2.Range("A1").Value = "P [MPa]"

$ws2.Range("A2").Value = 0.74
$ws2.Range("A3").Value = 1.32
$ws2.Range("A4").Value = 2.78
$ws2.Range("A5").Value = 3.31

$ws2.Range("B2").Value = 0.00022
$ws2.Range("B3").Value = 0.00038
$ws2.Range("B4").Value = 0.00073
$ws2.Range("B5").Value = 0.0009

$ws2.Range("B2:B5").NumberFormat = "0.00E+00"

$ws2.PageSetup.LeftMargin = 50.4
$ws2.PageSetup.RightMargin = 50.4
$ws2.PageSetup.TopMargin = 54
$ws2.PageSetup.BottomMargin = 54
$ws2.PageSetup.HeaderMargin = 21.6
$ws2.PageSetup.FooterMargin = 21.6

# ---------------------------------------------------------------
# Selections: CH4 keeps A2 selected (it is not the active tab),
# H2 becomes the active tab with B6 selected.
# ---------------------------------------------------------------
$ws1.Range("A2").Select() | Out-Null
$ws2.Select() | Out-Null
$ws2.Range("B6").Select() | Out-Null
